# Generate Report for Handoff
#
# Refresh the handoff identifiers (old GUID -> new GUID) and timestamps
# shown on the Overview / zh-cn / de-de sheets after a new handoff report
# was generated. The hyperlink targets themselves are untouched - only the
# visible text (cell value + hyperlink display text) is refreshed to show
# the new handoff file names and times.

$wb = $excel.ActiveWorkbook

$oldGuid = "2fc694ed-e560-416a-99f9-721d6e0b5ca7"
$newGuid = "f503f3e5-b93b-4bb6-bfa6-d43407f446b4"

$oldHash = "9723e6dc1d763b8e711ea5e7c678da2406a862f7"
$newHash = "ea1e835baed86b6be31222fd2d552adacb817c06"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet: A2 (handoff markdown link), D2 (handoff date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/3277ff6ded69c3e0e4c9ee968712c948c1b7eb0f/e2e/$oldGuid.md"

$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName) | Out-Null

$wsOverview.Range("D2").Value = "2016-52-20 22:52:35"

# ---------------------------------------------------------------------
# zh-cn sheet: A2 (md link), B2 (.md link), D2 (xlf link), E2 (xlf date)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f931fcefbc3859e9a6e86c9bb791ca43af7168b4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"

$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdAddress, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, "", "", $newZhXlfName) | Out-Null

$wsZh.Range("E2").Value = "2016-03-20 22:52:30"

# ---------------------------------------------------------------------
# de-de sheet: A2 (md link), B2 (.md link), D2 (xlf link), E2 (xlf date)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/da6317989b3988757d7618dcb15118b984fcf4dc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdAddress, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, "", "", $newDeXlfName) | Out-Null

$wsDe.Range("E2").Value = "2016-03-20 22:52:35"
